# Generate Report for Archive
#
# The localization status report moved on from the handoff stage, so every
# cell that still reads "Ready for handoff" is updated to "In Translation".
# Shrinking that text also lets Excel's column auto-fit narrow the Status
# columns that hold it (the "zh-cn"/"de-de" columns on the Overview sheet and
# the "Status" column on each per-language sheet).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    # Replace the old status text with the new one wherever it appears on
    # this sheet, and remember which columns actually held it so we can
    # re-fit their width afterwards (mirrors what Excel does automatically
    # when a user edits cell text that columns were auto-sized to).
    $usedRange = $ws.UsedRange
    $columnsToFit = @{}

    for ($r = 1; $r -le $usedRange.Rows.Count; $r++) {
        for ($c = 1; $c -le $usedRange.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Note: the status literal must be the left-hand operand of -eq;
            # some cells hold boolean values (e.g. "True"/"False" flags) and
            # PowerShell's -eq coerces the right-hand side to the left-hand
            # operand's type, so a boolean-left comparison would wrongly
            # treat every truthy cell as a match.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
                $columnsToFit[$c] = $true
            }
        }
    }

    foreach ($col in $columnsToFit.Keys) {
        $ws.Columns.Item($col).AutoFit() | Out-Null
    }
}
